# Burndown.xlsx - "adding new finished tasks"
# Add newly logged hours on the "Sprint 0" sheet:
#   - I6  (Banco de dados / day 7)      = 7
#   - H13 (Outros ... / day 6)          = 16
# Dependent SUM/running-total formulas (B6, B13, H14:W14, and the chart's
# cached series for 'Sprint 0'!C14:W14) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 0")

$ws.Range("I6").Value = 7
$ws.Range("H13").Value = 16

# Move the active selection, matching the author's last cursor position.
$ws.Activate() | Out-Null
$ws.Range("E25").Select() | Out-Null
